$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the inventory counts on row 2 (rounded totals)
$ws.Range("A2").Value = 992
$ws.Range("B2").Value = 952
$ws.Range("C2").Value = 952
$ws.Range("D2").Value = 952
$ws.Range("E2").Value = 985
$ws.Range("F2").Value = 986
